# Update odds/spread values per row based on the source diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.6
$ws.Range("H2").Value = 2.75
$ws.Range("I2").Value = 2.38
$ws.Range("J2").Value = 4.5
$ws.Range("K2").Value = 1.8
$ws.Range("L2").Value = 3.25
$ws.Range("M2").Value = 1.17
$ws.Range("N2").Value = 5
$ws.Range("O2").Value = 1.67
$ws.Range("P2").Value = 2.1
$ws.Range("Q2").Value = 2.32
$ws.Range("R2").Value = 1.62
$ws.Range("S2").Value = 3.4
$ws.Range("T2").Value = 1.33
$ws.Range("U2").Value = 5.6
$ws.Range("V2").Value = 1.15
$ws.Range("W2").Value = 7
$ws.Range("X2").Value = 1.1
$ws.Range("Y2").Value = 1.73
$ws.Range("Z2").Value = 2.08
$ws.Range("AA2").Value = 2.38
$ws.Range("AB2").Value = 1.53
$ws.Range("AC2").Value = 7.5
$ws.Range("AI2").Value = 5
$ws.Range("AJ2").Value = 5.5
$ws.Range("AL2").Value = 101
$ws.Range("AO2").Value = 9.5
$ws.Range("AP2").Value = 11
$ws.Range("AQ2").Value = 23
$ws.Range("AR2").Value = 26

# Row 3
$ws.Range("G3").Value = 1.73
$ws.Range("H3").Value = 3.3
$ws.Range("I3").Value = 5.5
$ws.Range("J3").Value = 2.5
$ws.Range("K3").Value = 1.95
$ws.Range("Q3").Value = 1.93
$ws.Range("R3").Value = 1.93
$ws.Range("AE3").Value = 9.5
$ws.Range("AF3").Value = 13
$ws.Range("AG3").Value = 19
$ws.Range("AJ3").Value = 6.5
$ws.Range("AL3").Value = 101
$ws.Range("AN3").Value = 10
$ws.Range("AS3").Value = 67

# Row 4
$ws.Range("G4").Value = 1.62
$ws.Range("H4").Value = 3.7
$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 2.3
$ws.Range("K4").Value = 2.05
$ws.Range("L4").Value = 6.5
$ws.Range("N4").Value = 7.5
$ws.Range("O4").Value = 1.4
$ws.Range("P4").Value = 2.75
$ws.Range("Q4").Value = 1.78
$ws.Range("R4").Value = 2.1
$ws.Range("S4").Value = 2.35
$ws.Range("T4").Value = 1.57
$ws.Range("U4").Value = 3.7
$ws.Range("V4").Value = 1.28
$ws.Range("W4").Value = 4.33
$ws.Range("X4").Value = 1.2
$ws.Range("Y4").Value = 1.5
$ws.Range("Z4").Value = 2.5
$ws.Range("AI4").Value = 7.5
$ws.Range("AK4").Value = 21
$ws.Range("AN4").Value = 12
$ws.Range("AO4").Value = 29

# Row 5
$ws.Range("S5").Value = 2.35
$ws.Range("T5").Value = 1.57

# Row 6
$ws.Range("H6").Value = 2.9
$ws.Range("I6").Value = 2.8
$ws.Range("J6").Value = 3.6
$ws.Range("AN6").Value = 7
$ws.Range("AO6").Value = 12

# Row 8
$ws.Range("G8").Value = 3.5
$ws.Range("I8").Value = 2.25
$ws.Range("J8").Value = 4.33
$ws.Range("L8").Value = 3.1
$ws.Range("M8").Value = 1.13
$ws.Range("N8").Value = 6
$ws.Range("O8").Value = 1.53
$ws.Range("P8").Value = 2.38
$ws.Range("Q8").Value = 2.03
$ws.Range("R8").Value = 1.83
$ws.Range("S8").Value = 2.7
$ws.Range("T8").Value = 1.44
$ws.Range("U8").Value = 4.3
$ws.Range("V8").Value = 1.21
$ws.Range("Y8").Value = 1.62
$ws.Range("Z8").Value = 2.2
$ws.Range("AD8").Value = 15
$ws.Range("AI8").Value = 6
$ws.Range("AO8").Value = 9.5
$ws.Range("AR8").Value = 23

# Row 13
$ws.Range("H13").Value = 3.4
$ws.Range("I13").Value = 4
$ws.Range("J13").Value = 2.63
$ws.Range("K13").Value = 2.05
$ws.Range("L13").Value = 4.75
$ws.Range("N13").Value = 8.5
$ws.Range("S13").Value = 2.2
$ws.Range("T13").Value = 1.67
$ws.Range("AA13").Value = 1.91
$ws.Range("AB13").Value = 1.8
$ws.Range("AC13").Value = 6.5
$ws.Range("AD13").Value = 8.5
$ws.Range("AE13").Value = 9
$ws.Range("AI13").Value = 8.5
$ws.Range("AK13").Value = 17
$ws.Range("AM13").Value = 351
$ws.Range("AN13").Value = 10

# Row 15
$ws.Range("W15").Value = 3.5
$ws.Range("X15").Value = 1.3
$ws.Range("AD15").Value = 10
$ws.Range("AS15").Value = 34

# Row 18
$ws.Range("G18").Value = 1.22
$ws.Range("M18").Value = 1.04
$ws.Range("N18").Value = 13

# Row 19
$ws.Range("I19").Value = 3.7
$ws.Range("AD19").Value = 9
$ws.Range("AF19").Value = 17
$ws.Range("AN19").Value = 10
$ws.Range("AO19").Value = 19
$ws.Range("AR19").Value = 34

# Row 24
$ws.Range("S24").Value = 2.05
$ws.Range("T24").Value = 1.75
